$wb = $excel.ActiveWorkbook

# --- Update the conversion text on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.71 = 6246.67 pesos`n✅ 6246.67 pesos = 1.7 = 953.07 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 585
$ws2.Range("O10").Value = 3654.3
$ws2.Range("N12").Value = 3677
$ws2.Range("O12").Value = 561.01
